$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Sony Headphones" -> "Sony", and update Count/Price/TotalCost
$ws.Range("A3").Value = "Sony"
$ws.Range("D3").Value = 25
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 272

# Row 4 used to hold a duplicate product row; now it only keeps the
# TotalCost subtotal in column F. Clear A4:E4 and set F4.
$ws.Range("A4:E4").ClearContents()
$ws.Range("F4").Value = 512

# Rows 5 and 6 (second duplicate row + grand total row) are removed entirely.
$ws.Range("A5:F6").EntireRow.Delete()
